$wb = $excel.ActiveWorkbook

# --- 1. Text change: "Ready for handoff" -> "In Translation" -------------
# This shared string is used by the "Status" columns on every sheet:
#   Overview!E2, Overview!F2, zh-cn!C2, de-de!C2
$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$overview.Range("E2").Value2 = "In Translation"
$overview.Range("F2").Value2 = "In Translation"
$zhcn.Range("C2").Value2     = "In Translation"
$dede.Range("C2").Value2     = "In Translation"

# --- 2. Column width changes (Status columns got narrower) ---------------
# Overview: columns E (5) and F (6)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de: column C (3)
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
